# Auto-generated edit script: updates cryptos list per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.735.54"
$ws.Range("E2").Value = "  -2.09%  "

$ws.Range("D3").Value = "3.494.80"
$ws.Range("E3").Value = "  -2.64%  "

$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").Value = "'570.16"
$ws.Range("E5").Value = "  -1.12%  "

$ws.Range("D6").Value = "'182.70"
$ws.Range("E6").Value = "  -3.50%  "

$ws.Range("D7").Value = "'0.613"
$ws.Range("E7").Value = "  -2.97%  "

$ws.Range("D8").Value = "3.488.69"
$ws.Range("E8").Value = "  -2.74%  "

$ws.Range("E9").Value = "  +0.13%  "

$ws.Range("D10").Value = "'0.184"
$ws.Range("E10").Value = "  +2.86%  "

$ws.Range("D11").Value = "'0.637"
$ws.Range("E11").Value = "  -3.67%  "

$ws.Range("D12").Value = "'53.51"
$ws.Range("E12").Value = "  -4.49%  "

$ws.Range("D13").Value = "'0.0000299"
$ws.Range("E13").Value = "  -0.92%  "

$ws.Range("D14").Value = "'9.40"
$ws.Range("E14").Value = "  -2.05%  "

$ws.Range("D15").Value = "4.066.86"
$ws.Range("E15").Value = "  -2.30%  "

$ws.Range("D16").Value = "'19.21"
$ws.Range("E16").Value = "  -3.64%  "

$ws.Range("D17").Value = "3.498.50"
$ws.Range("E17").Value = "  -2.39%  "

$ws.Range("D18").Value = "68.757.15"
$ws.Range("E18").Value = "  -1.89%  "

$ws.Range("D19").Value = "'12.44"
$ws.Range("E19").Value = "  -1.73%  "

$ws.Range("D20").Value = "'0.119"
$ws.Range("E20").Value = "  -1.30%  "

$ws.Range("D21").Value = "'538.75"
$ws.Range("E21").Value = "  +13.33%  "

$ws.Range("D22").Value = "'1.02"
$ws.Range("E22").Value = "  -2.03%  "

$ws.Range("D23").Value = "'20.11"
$ws.Range("E23").Value = "  +4.90%  "

$ws.Range("D24").Value = "'5.00"
$ws.Range("E24").Value = "  -1.76%  "

$ws.Range("D25").Value = "'4.39"
$ws.Range("E25").Value = "  +0.69%  "

$ws.Range("D26").Value = "'94.25"
$ws.Range("E26").Value = "  +6.14%  "

$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").Value = "'2.91"
$ws.Range("E27").Value = "  -4.41%  "

$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'10.86"
$ws.Range("E28").Value = "  -1.82%  "

$ws.Range("D29").Value = "'9.06"
$ws.Range("E29").Value = "  -3.35%  "

$ws.Range("D30").Value = "'31.34"
$ws.Range("E30").Value = "  -2.65%  "

$ws.Range("D31").Value = "'7.19"
$ws.Range("E31").Value = "  -6.75%  "

$ws.Range("D32").Value = "'12.56"
$ws.Range("E32").Value = "  +3.32%  "

$ws.Range("D33").Value = "'64.18"
$ws.Range("E33").Value = "  -3.02%  "

$ws.Range("E34").Value = "  -5.41%  "

$ws.Range("D35").Value = "'574.27"
$ws.Range("E35").Value = "  -1.72%  "

$ws.Range("D36").Value = "'37.81"
$ws.Range("E36").Value = "  -3.14%  "

$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").Value = "'0.998"
$ws.Range("E37").Value = "  -0.20%  "

$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").Value = "'0.397"
$ws.Range("E38").Value = "  -0.26%  "

$ws.Range("D39").Value = "'3.02"
$ws.Range("E39").Value = "  +5.34%  "

$ws.Range("D40").Value = "0.0₃0761"
$ws.Range("E40").Value = "  -4.80%  "

$ws.Range("D41").Value = "'0.133"
$ws.Range("E41").Value = "  -4.44%  "

$ws.Range("D42").Value = "'3.08"
$ws.Range("E42").Value = "  -4.44%  "

$ws.Range("D43").Value = "'3.32"
$ws.Range("E43").Value = "  -4.83%  "

$ws.Range("B44").Value = "ThetaToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D44").Value = "'2.96"
$ws.Range("E44").Value = "  -5.20%  "

$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "3.198.81"
$ws.Range("E45").Value = "  -1.38%  "

$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "'3.47"
$ws.Range("E46").Value = "  +4.33%  "

$ws.Range("D47").Value = "'0.0438"
$ws.Range("E47").Value = "  -0.77%  "

$ws.Range("D48").Value = "'9.08"
$ws.Range("E48").Value = "  -4.54%  "

$ws.Range("E49").Value = "  -2.57%  "

$ws.Range("E50").Value = "  -0.06%  "

$ws.Range("D51").Value = "'136.09"
$ws.Range("E51").Value = "  -0.66%  "

